# Trade #4 closed at 2026-02-16 21:20:43 - momentum DOWN +0.000%
#
# Inserts a new "momentum" worksheet between "leadlag" and "All Trades",
# populates it with the trade-log header row plus the new open trade
# (Trade #4), matching the column layout used by the other strategy
# sheets (e.g. "leadlag").

$wb = $excel.ActiveWorkbook

# Insert the new sheet right before "All Trades" so the final order is:
# Summary, leadlag, momentum, All Trades, Comparison
$beforeSheet = $wb.Worksheets.Item("All Trades")
$ws = $wb.Worksheets.Add($beforeSheet)
$ws.Name = "momentum"

# --- column widths (character units), matching the target layout -----
$widths = @(9, 12, 10, 10, 6, 13, 12, 8, 7, 7, 12, 44, 13, 16)
for ($i = 0; $i -lt $widths.Length; $i++) {
    $ws.Columns.Item($i + 1).ColumnWidth = ($widths[$i] - 0.8333333333333334)
}

# --- header row ---------------------------------------------------------
$headers = @("Trade #", "Date", "Time", "Strategy", "Side", "Entry Price", "Exit Price", "Status", "P&L %", "P&L $", "Confidence", "Entry Reason", "Exit Reason", "Duration (min)")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- trade #4 data row ---------------------------------------------------
$ws.Range("A2").Value = 4
# Force the date column to stay plain text ("2026-02-16") instead of
# being auto-parsed into a date serial number.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2026-02-16"
$ws.Range("C2").Value = "21:20:43"
$ws.Range("D2").Value = "momentum"
$ws.Range("E2").Value = "DOWN"
$ws.Range("F2").Value = 69460.925
$ws.Range("G2").Value = ""
$ws.Range("H2").Value = "OPEN"
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0.9
$ws.Range("L2").Value = "Downward momentum: -0.208% over 10 samples"
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = 0

# Restore the originally-active sheet (adding a sheet makes it active by
# default); the workbook was originally opened on the first sheet.
$wb.Worksheets.Item("Summary").Activate()
